$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("Q2").Value = 2.03
$ws.Range("R2").Value = 1.68
$ws.Range("AT2").Value = 2.62

# Row 3
$ws.Range("Q3").Value = 2.3
$ws.Range("R3").Value = 1.53

# Row 5
$ws.Range("M5").Value = 1.05
$ws.Range("N5").Value = 11
$ws.Range("O5").Value = 1.29
$ws.Range("P5").Value = 3.5
$ws.Range("Q5").Value = 1.92
$ws.Range("R5").Value = 1.82

# Row 6
$ws.Range("M6").Value = 1.02
$ws.Range("N6").Value = 12.8
$ws.Range("O6").Value = 1.15
$ws.Range("P6").Value = 4.1
$ws.Range("T6").Value = 3.28
$ws.Range("U6").Value = 1.84
$ws.Range("V6").Value = 1.92

# Row 7
$ws.Range("G7").Value = 2.62
$ws.Range("H7").Value = 3.1
$ws.Range("J7").Value = 3.15
$ws.Range("K7").Value = 2.07
$ws.Range("R7").Value = 1.93
$ws.Range("U7").Value = 1.52
$ws.Range("V7").Value = 2.2
$ws.Range("W7").Value = 10
$ws.Range("X7").Value = 15
$ws.Range("Z7").Value = 32
$ws.Range("AB7").Value = 24
$ws.Range("AC7").Value = 11
$ws.Range("AD7").Value = 6.2
$ws.Range("AE7").Value = 11
$ws.Range("AF7").Value = 40
$ws.Range("AG7").Value = 250
$ws.Range("AH7").Value = 10.25
$ws.Range("AN7").Value = 4.7
$ws.Range("AO7").Value = 14
$ws.Range("AQ7").Value = 60
$ws.Range("AT7").Value = 2.8
$ws.Range("AU7").Value = 6.2
$ws.Range("AV7").Value = 45
$ws.Range("BA7").Value = 75

# Row 9
$ws.Range("G9").Value = 1.53
$ws.Range("H9").Value = 4
$ws.Range("I9").Value = 6
$ws.Range("J9").Value = 2.05
$ws.Range("L9").Value = 5.5
$ws.Range("M9").Value = 1.04
$ws.Range("N9").Value = 13
$ws.Range("O9").Value = 1.2
$ws.Range("P9").Value = 4.33
$ws.Range("Q9").Value = 1.67
$ws.Range("R9").Value = 2.15
$ws.Range("U9").Value = 1.73
$ws.Range("V9").Value = 2
$ws.Range("W9").Value = 8
$ws.Range("X9").Value = 8
$ws.Range("Z9").Value = 11
$ws.Range("AD9").Value = 8
$ws.Range("AE9").Value = 15
$ws.Range("AG9").Value = 201
$ws.Range("AH9").Value = 19
$ws.Range("AI9").Value = 34
$ws.Range("AJ9").Value = 19
$ws.Range("AK9").Value = 67
$ws.Range("AL9").Value = 41
$ws.Range("AM9").Value = 41
$ws.Range("AN9").Value = 3.6
$ws.Range("AO9").Value = 7.5
$ws.Range("AQ9").Value = 21
$ws.Range("AS9").Value = 101
$ws.Range("AU9").Value = 8
$ws.Range("AW9").Value = 7.5
$ws.Range("AX9").Value = 29
$ws.Range("AZ9").Value = 101
$ws.Range("BA9").Value = 101
$ws.Range("BB9").Value = 201

# Row 12
$ws.Range("K12").Value = 1.92

# Row 13
$ws.Range("J13").Value = 1.87
$ws.Range("K13").Value = 2.37
$ws.Range("Q13").Value = 1.82
$ws.Range("R13").Value = 1.92
